$d = $word.ActiveDocument

function Split-RunAt($searchRange, $searchText, $charIndexInMatch) {
    # Finds $searchText within $searchRange (a Range), then narrows the
    # found range down to the single character at zero-based offset
    # $charIndexInMatch within the match, and toggles bold on/off on
    # that single character. Word's COM layer reacts to a real (even if
    # transient) formatting change by splitting the run(s) that cover
    # the target range into three runs: [before][char][after], each
    # carrying a copy of the original run properties.
    $rng = $searchRange.Duplicate
    $found = $rng.Find.Execute($searchText)
    if (-not $found) {
        throw "Could not find '$searchText'"
    }
    $matchLen = $rng.Text.Length
    # Move the start forward so only the target character (and anything
    # after it) remains, then pin the end back so only that one
    # character is left selected.
    [void]$rng.MoveStart(1, $charIndexInMatch)
    $rng.End = $rng.Start + 1
    $rng.Font.Bold = $true
    $rng.Font.Bold = $false
}

# --- 1) "algorithms - Combining merge sort and insertion sort - Computer
#        Science Stack Exchange" hyperlink display text (field result) ---
# First normalize the run via a full-text Find & Replace (this clears any
# leftover rsid metadata on the run), then split off the single "t" that
# turns "sor" into "sort".
$p1 = $d.Paragraphs(1).Range
[void]$p1.Find.Execute(
    "algorithms - Combining merge sort and insertion sort - Computer Science Stack Exchange",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "algorithms - Combining merge sort and insertion sort - Computer Science Stack Exchange",
    2)

Split-RunAt $d.Paragraphs(1).Range "merge sort" 9

# --- 2) "Algorithms: GATE CSE 2007 | Question: 41 (gateoverflow.in)"
#        hyperlink display text ---
$p5 = $d.Paragraphs(5).Range
[void]$p5.Find.Execute(
    "Algorithms: GATE CSE 2007 | Question: 41 (gateoverflow.in)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Algorithms: GATE CSE 2007 | Question: 41 (gateoverflow.in)",
    2)

Split-RunAt $d.Paragraphs(5).Range "gateoverflow" 6

Write-Host "Done: split hyperlink display-text runs in paragraphs 1 and 5."
